# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" block together
# with its surrounding blank paragraph and the following page-break
# paragraph, which sat right after the "Requisitos" section's
# "LOT2004: Bioquímica (Requisito fraco)" paragraph.
#
# Before:
#   ... LOT2004: Bioquímica (Requisito fraco)
#   <empty>
#   Ver no Jupiter Salvar em pdf Salvar em docx
#   <empty>
#   <empty, pageBreakBefore, jc=left>
#   <empty>
#   <empty, pageBreakBefore>
#
# After:
#   ... LOT2004: Bioquímica (Requisito fraco)
#   <empty>
#   <empty, pageBreakBefore>

$d = $word.ActiveDocument

$target = "Ver no Jupiter Salvar em pdf Salvar em docx"

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "$target*") {
        # Delete this paragraph plus the blank paragraph immediately
        # before it and the two paragraphs immediately after it (a blank
        # one and the page-break one), leaving the rest untouched.
        $startPara = $i - 1
        $endPara = $i + 2

        $start = $d.Paragraphs($startPara).Range.Start
        $end = $d.Paragraphs($endPara).Range.End

        $r = $d.Range($start, $end)
        $r.Delete()

        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the '$target' paragraph"
}
